$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '37.364.51'
$ws.Range('E2').Value = '  +2.25%  '

$ws.Range('D3').Value = '2.084.08'
$ws.Range('E3').Value = '  -0.66%  '

$ws.Range('E4').Value = '  +0.05%  '

$ws.Range('D5').NumberFormat = '@'
$ws.Range('D5').Value = '251.46'
$ws.Range('E5').Value = '  +1.22%  '

$ws.Range('D6').NumberFormat = '@'
$ws.Range('D6').Value = '0.665'
$ws.Range('E6').Value = '  -0.60%  '

$ws.Range('E7').Value = '  +0.00%  '

$ws.Range('D8').NumberFormat = '@'
$ws.Range('D8').Value = '57.38'
$ws.Range('E8').Value = '  +26.87%  '

$ws.Range('D9').NumberFormat = '@'
$ws.Range('D9').Value = '62.47'
$ws.Range('E9').Value = '  +1.91%  '

$ws.Range('D10').NumberFormat = '@'
$ws.Range('D10').Value = '0.383'
$ws.Range('E10').Value = '  +4.79%  '

$ws.Range('D11').NumberFormat = '@'
$ws.Range('D11').Value = '0.0753'
$ws.Range('E11').Value = '  +3.16%  '

$ws.Range('D12').NumberFormat = '@'
$ws.Range('D12').Value = '0.107'
$ws.Range('E12').Value = '  +7.43%  '

$ws.Range('D13').NumberFormat = '@'
$ws.Range('D13').Value = '15.39'
$ws.Range('E13').Value = '  +5.64%  '

$ws.Range('D14').Value = '2.393.67'
$ws.Range('E14').Value = '  -0.44%  '

$ws.Range('D15').NumberFormat = '@'
$ws.Range('D15').Value = '0.850'
$ws.Range('E15').Value = '  +1.30%  '

$ws.Range('D16').NumberFormat = '@'
$ws.Range('D16').Value = '5.28'
$ws.Range('E16').Value = '  +5.18%  '

$ws.Range('D17').Value = '2.093.93'
$ws.Range('E17').Value = '  -0.08%  '

$ws.Range('D18').Value = '37.332.39'
$ws.Range('E18').Value = '  +1.94%  '

$ws.Range('D19').NumberFormat = '@'
$ws.Range('D19').Value = '73.24'
$ws.Range('E19').Value = '  +0.93%  '

$ws.Range('D20').NumberFormat = '@'
$ws.Range('D20').Value = '14.79'
$ws.Range('E20').Value = '  +15.21%  '

$ws.Range('D21').Value = '0.0₃0851'
$ws.Range('E21').Value = '  +4.00%  '

$ws.Range('D22').NumberFormat = '@'
$ws.Range('D22').Value = '240.90'
$ws.Range('E22').Value = '  +0.05%  '

$ws.Range('D23').NumberFormat = '@'
$ws.Range('D23').Value = '5.27'
$ws.Range('E23').Value = '  +4.79%  '

$ws.Range('D24').NumberFormat = '@'
$ws.Range('D24').Value = '0.999'
$ws.Range('E24').Value = '  -0.12%  '

$ws.Range('E25').Value = '  +0.97%  '

$ws.Range('D26').NumberFormat = '@'
$ws.Range('D26').Value = '171.61'
$ws.Range('E26').Value = '  +1.02%  '

$ws.Range('D27').NumberFormat = '@'
$ws.Range('D27').Value = '9.27'
$ws.Range('E27').Value = '  +3.90%  '

$ws.Range('D28').NumberFormat = '@'
$ws.Range('D28').Value = '21.10'
$ws.Range('E28').Value = '  +2.79%  '

$ws.Range('E29').Value = '  +1.16%  '

$ws.Range('E30').Value = '  +1.62%  '

$ws.Range('B31').Value = 'ImmutableX'
$ws.Range('C31').Value = 'https://coinranking.com/coin/Z96jIvLU7+immutablex-imx'
$ws.Range('D31').NumberFormat = '@'
$ws.Range('D31').Value = '1.11'
$ws.Range('E31').Value = '  +22.53%  '

$ws.Range('B32').Value = 'Gas'
$ws.Range('C32').Value = 'https://coinranking.com/coin/hfw0nnnLtSFc7+gas-gas'
$ws.Range('D32').NumberFormat = '@'
$ws.Range('D32').Value = '23.01'
$ws.Range('E32').Value = '  +2.81%  '

$ws.Range('D33').NumberFormat = '@'
$ws.Range('D33').Value = '4.58'
$ws.Range('E33').Value = '  +3.45%  '

$ws.Range('D34').NumberFormat = '@'
$ws.Range('D34').Value = '0.0628'
$ws.Range('E34').Value = '  +5.81%  '

$ws.Range('D35').NumberFormat = '@'
$ws.Range('D35').Value = '0.0911'
$ws.Range('E35').Value = '  +0.10%  '

$ws.Range('D36').NumberFormat = '@'
$ws.Range('D36').Value = '4.32'
$ws.Range('E36').Value = '  +6.51%  '

$ws.Range('E37').Value = '  +0.05%  '

$ws.Range('D38').NumberFormat = '@'
$ws.Range('D38').Value = '2.29'
$ws.Range('E38').Value = '  -0.75%  '

$ws.Range('D39').NumberFormat = '@'
$ws.Range('D39').Value = '1.85'
$ws.Range('E39').Value = '  -1.44%  '

$ws.Range('D40').NumberFormat = '@'
$ws.Range('D40').Value = '1.35'
$ws.Range('E40').Value = '  -0.14%  '

$ws.Range('E41').Value = '  +5.05%  '

$ws.Range('B42').Value = 'Cronos'
$ws.Range('C42').Value = 'https://coinranking.com/coin/65PHZTpmE55b+cronos-cro'
$ws.Range('D42').NumberFormat = '@'
$ws.Range('D42').Value = '0.0993'
$ws.Range('E42').Value = '  +19.19%  '

$ws.Range('B43').Value = 'InjectiveProtocol'
$ws.Range('C43').Value = 'https://coinranking.com/coin/PkY9BmsyW+injectiveprotocol-inj'
$ws.Range('D43').NumberFormat = '@'
$ws.Range('D43').Value = '17.79'
$ws.Range('E43').Value = '  +10.04%  '

$ws.Range('E44').Value = '  -1.23%  '

$ws.Range('D45').NumberFormat = '@'
$ws.Range('D45').Value = '100.20'
$ws.Range('E45').Value = '  +1.04%  '

$ws.Range('D46').NumberFormat = '@'
$ws.Range('D46').Value = '4.46'
$ws.Range('E46').Value = '  +114.63%  '

$ws.Range('D47').NumberFormat = '@'
$ws.Range('D47').Value = '2.81'
$ws.Range('E47').Value = '  +0.19%  '

$ws.Range('D48').Value = '1.328.48'
$ws.Range('E48').Value = '  -2.36%  '

$ws.Range('D49').NumberFormat = '@'
$ws.Range('D49').Value = '2.95'
$ws.Range('E49').Value = '  +4.27%  '

$ws.Range('D50').NumberFormat = '@'
$ws.Range('D50').Value = '2.37'
$ws.Range('E50').Value = '  +4.97%  '

$ws.Range('D51').NumberFormat = '@'
$ws.Range('D51').Value = '7.05'
$ws.Range('E51').Value = '  +9.78%  '
